# This script applies a cyclic shift to rows 2-6 on the active worksheet:
# the data previously in row 6 moves to row 2, and the data previously in
# rows 2-5 each shift down one row (to rows 3-6 respectively). Only the
# columns A, B, D, E, F, G, H, Q, R participate (the rest are identical
# across these rows already).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that change, for rows 2..6.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$original = @{}
foreach ($r in 2..6) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowVals
}

# New row r (3..6) gets the old values of row r-1; new row 2 gets the old
# values of row 6 (cyclic shift down).
foreach ($r in 3..6) {
    $src = $original[$r - 1]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $src[$c]
    }
}
foreach ($c in $cols) {
    $ws.Range("${c}2").Value = $original[6][$c]
}
